# Update the "Forecast Comparison" sheet:
#  - insert a new "Week_Start_Date" column between "Week" and "ASIN"
#  - shorten the week labels ("W01" -> "W1", ... "W09" -> "W9")
#  - correct a handful of MyForecast values
#  - the trailing "is_holiday_week" column becomes boolean (already implied
#    by assigning $true/$false below)
# Then update the "Summary" sheet totals that depend on the corrected
# forecast numbers.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("Forecast Comparison")
$sum = $wb.Worksheets.Item("Summary")

# --- 1. Insert the new column B ("Week_Start_Date") ---------------------
$ws.Columns.Item(2).Insert()
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# --- 2. Fill in the week-start dates (kept as text, not Excel dates) ----
$weekStarts = @(
    "2025-01-05", "2025-01-12", "2025-01-19", "2025-01-26",
    "2025-02-02", "2025-02-09", "2025-02-16", "2025-02-23",
    "2025-03-02", "2025-03-09", "2025-03-16", "2025-03-23",
    "2025-03-30", "2025-04-06", "2025-04-13", "2025-04-20"
)
for ($i = 0; $i -lt $weekStarts.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = "'" + $weekStarts[$i]
}

# --- 3. Shorten the week labels in column A (W01 -> W1, ... W09 -> W9) --
for ($i = 1; $i -le 9; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = "W" + $i
}

# --- 4. Correct a few MyForecast values (now column D) -----------------
$ws.Cells.Item(6, 4).Value  = 30   # W5  : 31 -> 30
$ws.Cells.Item(9, 4).Value  = 30   # W8  : 29 -> 30
$ws.Cells.Item(13, 4).Value = 23   # W12 : 22 -> 23
$ws.Cells.Item(14, 4).Value = 23   # W13 : 22 -> 23

# --- 5. Make the is_holiday_week column (now J) boolean -----------------
for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 10).Value = $false
}

# --- 6. Update the dependent Summary totals -----------------------------
$sum.Cells.Item(9, 2).Value  = "'400"   # Total Forecast (16 Weeks): 399 -> 400
$sum.Cells.Item(10, 2).Value = "'209"   # Total Forecast (8 Weeks):  210 -> 209
$sum.Cells.Item(11, 2).Value = "'99"    # Total Forecast (4 Weeks):  100 -> 99
$sum.Cells.Item(12, 2).Value = "'30"    # Max Forecast:               31 -> 30
